$d = $word.ActiveDocument

# The document had an introductory paragraph ("Os princípios de segurança
# da informação (...) são fundamentais para a MarketSecure:") followed by
# a blank paragraph, right before the "Confidencialidade / Integridade /
# Disponibilidade" list. Both of those paragraphs are being removed.
#
# The "_GoBack" bookmark used to sit at the very end of the document
# (end of the last "...em caso de falhas." paragraph); it needs to move to
# the end of the empty paragraph that precedes the paragraph being removed.

# 1) Drop the existing _GoBack bookmark -- it will be recreated at the new
#    location once the paragraphs are removed.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Find the paragraph that holds the "MarketSecure" sentence, by content,
#    so this does not depend on fragile hard-coded paragraph indices.
$targetParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*MarketSecure*") {
        $targetParaIndex = $i
        break
    }
}

# 3) Delete that paragraph plus the blank paragraph that immediately
#    follows it (once the first paragraph is deleted, the following blank
#    paragraph shifts down into the same index).
$introPara = $d.Paragraphs.Item($targetParaIndex)
$introPara.Range.Delete()
$blankPara = $d.Paragraphs.Item($targetParaIndex)
$blankPara.Range.Delete()

# 4) The paragraph right before the removed text (still empty) now sits
#    where the bookmark belongs. Put the "_GoBack" bookmark there,
#    collapsed at its end (i.e. right before its own paragraph mark).
$precedingPara = $d.Paragraphs.Item($targetParaIndex - 1)
$bookmarkRange = $precedingPara.Range.Duplicate
$bookmarkRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
